# Edit "leaked datasets" sheet: extend the small reference table with more
# rows of breached-service data, add a "Year" column and a source link, and
# apply a bold header row with a double bottom border.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leaked datasets")

# --- Cell values -----------------------------------------------------------
# NB: values are written in this particular order so that new shared-string
# table entries come out in the same sequence as in the target workbook.

$ws.Range("A5").Value  = "MySpace"
$ws.Range("C1").Value  = "Year"
$ws.Range("A13").Value = "http://breachlevelindex.com/top-data-breaches"
$ws.Range("A6").Value  = "ebay"
$ws.Range("B6").Value  = "145 M"
$ws.Range("B5").Value  = "360 M"
$ws.Range("A7").Value  = "Adobe"
$ws.Range("A8").Value  = "Yahoo"
$ws.Range("B8").Value  = "1 B"
$ws.Range("B4").Value  = "68 M"
$ws.Range("B7").Value  = "36 M"
$ws.Range("B3").Value  = "164 M"
$ws.Range("A4").Value  = "Dropbox"

# Year numbers for each data source
$ws.Range("C2").Value = 2009
$ws.Range("C3").Value = 2016
$ws.Range("C4").Value = 2012
$ws.Range("C5").Value = 2013
$ws.Range("C6").Value = 2014
$ws.Range("C7").Value = 2013
$ws.Range("C8").Value = 2013

# --- Header row formatting --------------------------------------------------
$header = $ws.Range("A1:D1")
$header.Font.Bold = $true
$header.Borders.Item(9).LineStyle = -4119

# --- View settings -----------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Zoom = 160
$ws.Activate()
$ws.Range("A1:D8").Select()

Write-Output "Edit complete"
